$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (e.g. "51.90", "1.001")
# instead of being auto-converted to numbers by Excel, by forcing Text format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.822.72'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.740.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.54'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.38%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5162'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2722'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '38.74'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06085'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.743.03'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07005'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.15'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6298'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.494'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.25'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.842.55'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.43'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006605'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.961.79'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.058'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.408'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.082'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.502'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.815'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.95'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.60'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08305'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.617'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.373'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04398'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.607'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9665'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5952'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.677'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01553'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.933'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9995'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.33'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3797'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7244'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.873'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05482'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.182'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.15%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.76'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.90'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.42%  '
